$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the "Tanggal Lahir" example cell: date format sample changes ---
$ws.Range("F2").Value = "Ex : 17/02/2002"

# --- 2. New column U: "Tanggal Mulai Bekerja" (header) / example date ---
$ws.Range("U1").Value = "Tanggal Mulai Bekerja"
# Give U1 the same look as the other account-header cells (R1:T1): yellow
# fill + full thin border. Copy/paste-format from R1 is the reliable way to
# reuse that existing style instead of re-building it from scratch.
$ws.Range("R1").Copy()
$ws.Range("U1").PasteSpecial(-4122) | Out-Null

$ws.Range("U2").Value = "Ex : 01/01/2024"
$ws.Range("U2").NumberFormat = "mm-dd-yy"

# --- 3. New column V: "Nama Posisi" (header) / helper note ---
$ws.Range("V1").Value = "Nama Posisi"
# Start from the same bordered header style as U1/R1, then re-tint the fill
# to the new accent1 blue used for this column.
$ws.Range("R1").Copy()
$ws.Range("V1").PasteSpecial(-4122) | Out-Null
$ws.Range("V1").Interior.ThemeColor = 5
$ws.Range("V1").Interior.TintAndShade = 0.59999389629810485

$ws.Range("V2").Value = '(Wajib sama dengan data dari database, boleh lebih dari 1 dengan separator koma " , " )'

$excel.CutCopyMode = 0

# --- 4. Column widths for the two new columns ---
$ws.Columns.Item(21).ColumnWidth = 30.6
$ws.Columns.Item(22).ColumnWidth = 85.1

# --- 5. View state: scroll over and select U6, matching the saved selection ---
$excel.ActiveWindow.ScrollColumn = 16
$ws.Range("U6").Select()
